# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 273 (shifting the existing
# historical records at rows 273:385 down to rows 274:386).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 273. Excel shifts all data
# (and formatting) for rows 273:385 down to rows 274:386 automatically.
$ws.Rows.Item(273).Insert()

# Populate the newly inserted row 273 with the new weekly record.
$ws.Range("A273").Value = 10
$ws.Range("B273").Value = "Vega Modelo de Temuco"
$ws.Range("C273").Value = "La Araucanía"
$ws.Range("D273").Value = 44755
$ws.Range("D273").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E273").Value = 9
$ws.Range("F273").Value = "Fruta"
$ws.Range("G273").Value = 100108
$ws.Range("H273").Value = "Tropicales y subtropicales"
$ws.Range("I273").Value = 100108002
$ws.Range("J273").Value = "Mango"
$ws.Range("K273").Value = "Sin especificar"
$ws.Range("L273").Value = "Primera"
$ws.Range("M273").Value = 125
$ws.Range("N273").Value = 8000
$ws.Range("O273").Value = 8000
$ws.Range("P273").Value = 8000
$ws.Range("Q273").Value = "$/bandeja 4 kilos"
$ws.Range("R273").Value = "Brasil"
$ws.Range("S273").Value = 2000
$ws.Range("T273").Value = 4
